{"js": "// Into the Breach date update: change the series date from October to\n// January (and the weekday from Tuesday to Wednesday) for the first\n// session, and drop the now-obsolete per-topic dates from the schedule\n// paragraph later in the flyer.\n\nconst body = context.document.body;\n\n// 1) \"Tues 19 Oct 7 PM\" -> \"Wed 19 Jan 7 PM\"\nconst firstSessionHits = body.search(\"Tues 19 Oct 7 PM\", { matchCase: true });\nfirstSessionHits.load(\"items\");\nawait context.sync();\nfor (const hit of firstSessionHits.items) {\n  hit.insertText(\"Wed 19 Jan 7 PM\", \"Replace\");\n}\nawait context.sync();\n\n// 2) \"Enroll by 18 Oct:\" -> \"Enroll by 18 Jan:\"\nconst enrollHits = body.search(\"Enroll by 18 Oct:\", { matchCase: true });\nenrollHits.load(\"items\");\nawait context.sync();\nfor (const hit of enrollHits.items) {\n  hit.insertText(\"Enroll by 18 Jan:\", \"Replace\");\n}\nawait context.sync();\n\n// 3) Remove the per-topic parenthetical dates from the schedule sentence.\nconst scheduleHits = body.search(\n  \"Brotherhood (16 Nov), Leadership (18 Jan), Fatherhood (15 Feb), Family Life (15 Mar) and Prayer (19 Apr).\",\n  { matchCase: true }\n);\nscheduleHits.load(\"items\");\nawait context.sync();\nfor (const hit of scheduleHits.items) {\n  hit.insertText(\"Brotherhood, Leadership, Fatherhood, Family Life and Prayer.\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Into the Breach date update: change the series date from October to\n# January (and the weekday from Tuesday to Wednesday) for the first\n# session, and drop the now-obsolete per-topic dates from the schedule\n# paragraph later in the flyer.\n\n$d = $word.ActiveDocument\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n#   ReplaceWith, Replace)\n# Wrap: 1 = wdFindContinue ; Replace: 2 = wdReplaceAll\n\n# 1) \"Tues 19 Oct 7 PM\" -> \"Wed 19 Jan 7 PM\"\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Execute(\"Tues 19 Oct 7 PM\", $false, $false, $false, $false, $false, $true, 1, $false, \"Wed 19 Jan 7 PM\", 2)\n\n# 2) \"Enroll by 18 Oct:\" -> \"Enroll by 18 Jan:\"\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Execute(\"Enroll by 18 Oct:\", $false, $false, $false, $false, $false, $true, 1, $false, \"Enroll by 18 Jan:\", 2)\n\n# 3) Remove the per-topic parenthetical dates from the schedule sentence.\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Execute(\"Brotherhood (16 Nov), Leadership (18 Jan), Fatherhood (15 Feb), Family Life (15 Mar) and Prayer (19 Apr).\", $false, $false, $false, $false, $false, $true, 1, $false, \"Brotherhood, Leadership, Fatherhood, Family Life and Prayer.\", 2)\n"}
